$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 545; existing rows 545-609 shift down to 546-610.
$ws.Rows(545).Insert()

# Populate the newly inserted row 545 with the new "Mango" price record.
$ws.Cells.Item(545,1).Value  = 10
$ws.Cells.Item(545,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(545,3).Value  = "La Araucanía"
$ws.Cells.Item(545,4).Value  = 45124
$ws.Cells.Item(545,5).Value  = 9
$ws.Cells.Item(545,6).Value  = "Fruta"
$ws.Cells.Item(545,7).Value  = 100108
$ws.Cells.Item(545,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(545,9).Value  = 100108002
$ws.Cells.Item(545,10).Value = "Mango"
$ws.Cells.Item(545,11).Value = "Sin especificar"
$ws.Cells.Item(545,12).Value = "Primera"
$ws.Cells.Item(545,13).Value = 185
$ws.Cells.Item(545,14).Value = 9000
$ws.Cells.Item(545,15).Value = 9000
$ws.Cells.Item(545,16).Value = 9000
$ws.Cells.Item(545,17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(545,18).Value = "Brasil"
$ws.Cells.Item(545,19).Value = 2250
$ws.Cells.Item(545,20).Value = 4
